$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "c48e5a8a-7942-41de-a522-acb950ba8eda"
$ws.Range("D3").Value = "8e38ee81-1e04-49f1-b13b-f3133ac65281"
$ws.Range("D4").Value = "d0461977-e6c4-4add-a214-6e9fd42f285e"

$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMjI3NTk0MX0.xTUT3jJCjKmOvYExR3RG_X46xMhz9CZ1KmRKcug5M9I"
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMjI3NTk0Mn0.wDvr_CBYZqirz7frd6HFJZtSzqeVz984s08ZlgQbJ8g"
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDIyNzU5NDN9.jcUjSZSokv_DwQmxQIrdE0eChO6IfYliASffESHWq50"
